$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C12 from a date value to the text "DecV1.5"
$ws.Range("C12").ClearFormats()
$ws.Range("C12").Value = "DecV1.5"

# Update D12 from "60k" to "60K"
$ws.Range("D12").Value = "60K"

# Update the saved selection/active cell to J13
$ws.Range("J13").Select()
